# "Updated parsing of CPR Stockpile"
# Adds a new "Huadian Caofeidian" data column to the CPR Stockpile sheet
# (inserted just before the existing "CCI 4700" / "FOB Vostochny 5500"
# columns), fixes up the page setup, and leaves the CPR Stockpile sheet
# as the active tab/selection.

$wb = $excel.ActiveWorkbook

$wsStockpile = $wb.Worksheets.Item("CPR Stockpile")
$wsFreight   = $wb.Worksheets.Item("Freight")
$wsWeather   = $wb.Worksheets.Item("China Weather")

# --- Insert the new "Huadian Caofeidian" column at column I (9) ---------
$wsStockpile.Columns.Item(9).Insert() | Out-Null

$newHeader = $wsStockpile.Cells.Item(1, 9)
$newHeader.Value = "Huadian Caofeidian"

# Match formatting used by the other whole-number stockpile headers
# (bold slightly-tinted header font, grey fill, centered) — reuse the
# style already present on the "Freight" header row, then force the
# integer "#,##0" number format for this column.
$wsFreight.Cells.Item(1, 2).Copy() | Out-Null
$newHeader.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$newHeader.NumberFormat = "#,##0"

# The trailing "Update?" column (now column L) should match the style
# already used for that header on the "China Weather" sheet.
$wsWeather.Cells.Item(1, 6).Copy() | Out-Null
$wsStockpile.Cells.Item(1, 12).PasteSpecial(-4122) | Out-Null # xlPasteFormats

$excel.CutCopyMode = 0

# --- Page setup for CPR Stockpile ---------------------------------------
$wsStockpile.PageSetup.PaperSize = 9     # xlPaperA4
$wsStockpile.PageSetup.Orientation = 1   # xlPortrait

# --- Make CPR Stockpile the active sheet/selection ----------------------
$wsStockpile.Activate() | Out-Null
$wsStockpile.Range("D7").Select() | Out-Null
